$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the tnrsdate column (T) from 45905 to 45909 for rows 2 through 11
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Range("T$row")
    if ($cell.Value2 -eq 45905) {
        $cell.Value = 45909
    }
}
